# Rename the year-dependent generator-cost headers in row 1 from the
# "...Gyear20XX" form to the unified "...Gy20XX" form.
# Columns F (6) through AG (33) hold these 28 headers:
#   CostCapGyear / CostOperationVarGyear / CostOperationFixGyear / LifetimeGyear
# each for years 2020,2025,2030,2035,2040,2045,2050

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$prefixes = @("CostCapGy", "CostOperationVarGy", "CostOperationFixGy", "LifetimeGy")
$years = @(2020, 2025, 2030, 2035, 2040, 2045, 2050)

$col = 6  # column F
foreach ($prefix in $prefixes) {
    foreach ($year in $years) {
        $ws.Cells.Item(1, $col).Value = "$prefix$year"
        $col = $col + 1
    }
}

# Update the view state to match: scrolled so column W is left-most visible,
# with AH1 as the active/selected cell.
$ws.Activate()
$ws.Range("AH1").Select()
$excel.ActiveWindow.ScrollColumn = 23
$excel.ActiveWindow.ScrollRow = 1
